$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update the BMG note (cell I4) with the expanded note text, with rich formatting ---
$noteText = "From BMG Extruder`n*NOTE*`n1 - Some BMG clone have a 4mm thick gear, they will not fit.`n  Verified sources are: `n  a) Original BMG from BondTech`n  b) https://s.click.aliexpress.com/e/_Ao9yaZ`n2 -  Short side of the shaft might need to be filed down to avoid grinding against stepper motor."
$cellI4 = $ws.Range("I4")
$cellI4.Value = $noteText

$boldStart = $noteText.IndexOf("*NOTE*") + 1
$boldLen = $noteText.Length - $boldStart + 1
$boldChars = $cellI4.Characters($boldStart, $boldLen)
$boldChars.Font.Bold = $true
$boldChars.Font.Color = 255
$boldChars.Font.Name = "Calibri"
$boldChars.Font.Size = 11

# --- 2. Add Vendor (J) and Vendor URL (K) info for rows 2,3,4,5,6,7,9 ---
$vendorName = "Triangle Lab"
$vendorUrl = "https://s.click.aliexpress.com/e/_Ao9yaZ"

$ws.Range("K2:K3").Value = $vendorUrl
$ws.Range("K4").Value = $vendorUrl
$ws.Range("K5").Value = $vendorUrl
$ws.Range("K6:K7").Value = $vendorUrl
$ws.Range("K9").Value = $vendorUrl

$ws.Range("J2").Value = $vendorName
$ws.Range("J3").Value = $vendorName
$ws.Range("J4").Value = $vendorName
$ws.Range("J5").Value = $vendorName
$ws.Range("J6").Value = $vendorName
$ws.Range("J7").Value = $vendorName
$ws.Range("J9").Value = $vendorName

# Add hyperlinks
$ws.Hyperlinks.Add($ws.Range("K2:K3"), $vendorUrl, "", "", $vendorUrl)
$ws.Hyperlinks.Add($ws.Range("K4"), $vendorUrl)
$ws.Hyperlinks.Add($ws.Range("K5"), $vendorUrl)
$ws.Hyperlinks.Add($ws.Range("K6:K7"), $vendorUrl, "", "", $vendorUrl)
$ws.Hyperlinks.Add($ws.Range("K9"), $vendorUrl)

# Match formatting (font/border/style) of the new hyperlink cells to the existing hyperlink cell K13
$ws.Range("K13").Copy()
$ws.Range("K2:K7").PasteSpecial(-4122)
$ws.Range("K9").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- 3. Update view state: zoom and selection ---
$excel.ActiveWindow.Zoom = 55
$ws.Range("J2:K9").Select()

Write-Host "Edit complete"
